$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activity Log")

$ws.Range("A38").Value = "vaishali.kh2310@gmail.com"
$ws.Range("B38").Value = "Logout"
$ws.Range("C38").Value = "2025-07-04 15:44:25"

$ws.Range("A39").Value = "vaishali.kh2310@gmail.com"
$ws.Range("B39").Value = "Login"
$ws.Range("C39").Value = "2025-07-04 15:45:11"

$ws.Range("A40").Value = "vaishali.kh2310@gmail.com"
$ws.Range("B40").Value = "Logout"
$ws.Range("C40").Value = "2025-07-04 15:49:27"
